$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the attendance row for the new 10/6 / 4:15 meeting (row 19)
$ws.Range("B19").Value = "10/6 / 4:15"
$ws.Range("C19").Value = "Google Hangout"
$ws.Range("D19:I19").Value = "A"

# Match the cell formatting used by the other filled-in rows (copy the
# "top of block" border style from the row above, like row 18 already has)
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active cell selection
$ws.Range("J12").Select()
